# Insert a new price-record row at row 721 (pushing the existing rows 721:767
# down to 722:768) and populate it with the new "Ají" record:
#   Fecha=44714 (serial date), Variedad=Inferno, Volumen=65,
#   Precio min/max/prom=28000, Unidad=$/caja 15 kilos,
#   Origen=Región de Arica y Parinacota, Precio $/Kg=1867, Kg o Unidades=15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(721).Insert()

$ws.Cells.Item(721, 1).Value  = 10
$ws.Cells.Item(721, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(721, 3).Value  = "La Araucanía"
$ws.Cells.Item(721, 4).Value  = 44714
$ws.Cells.Item(721, 5).Value  = 9
$ws.Cells.Item(721, 6).Value  = 100112021
$ws.Cells.Item(721, 7).Value  = "Ají"
$ws.Cells.Item(721, 8).Value  = "Inferno"
$ws.Cells.Item(721, 9).Value  = "Primera"
$ws.Cells.Item(721, 10).Value = 65
$ws.Cells.Item(721, 11).Value = 28000
$ws.Cells.Item(721, 12).Value = 28000
$ws.Cells.Item(721, 13).Value = 28000
$ws.Cells.Item(721, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(721, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(721, 16).Value = 1867
$ws.Cells.Item(721, 17).Value = 15
$ws.Cells.Item(721, 18).Value = "Hortaliza"
